# Apply the "Updated vs latest scrape" edit to the CloudStore analysis workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the previously-missing "data" search numbers for Mar 2014 ---
# D3 / D4 used to hold "?" placeholders; now they have real counts, and the
# dependent ratio/variation formulas in E3/F3 stop erroring out.
$ws.Range("D3").Value = 3116
$ws.Range("D4").Value = 656
$ws.Range("E3").Formula = "=D3/D2"
$ws.Range("F3").Formula = "=D3/B3-1"

# --- 2. Update the A3 label: drop the "(1)" footnote marker ---
$cell = $ws.Range("A3")
$cell.Value = 'No. of products or services containing the word "data" in the description'
$text = $cell.Value2
$boldStart = $text.IndexOf('data') + 1
$boldLen = 4
$cell.Characters($boldStart, $boldLen).Font.Bold = $true
$suffixStart = $boldStart + $boldLen
$suffixLen = $text.Length - $suffixStart + 1
$cell.Characters($suffixStart, $suffixLen).Font.Size = 12

# --- 3. Remove the footnote row (old row 8, merged A8:F8) and the blank
#        spacer row beneath it (old row 9); give the remaining small spacer
#        row (row 7) the taller, formatted look the footnote row used to have ---
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

$ws.Rows.Item(7).RowHeight = 28
$ws.Range("A7:F7").Font.Size = 10
$ws.Range("A7:F7").Font.Name = "Calibri"
$ws.Range("A7:F7").VerticalAlignment = -4160
$ws.Range("A7:F7").WrapText = $true

# --- 4. Selection cosmetics (matches the saved selection in the workbook) ---
$ws.Range("A13").Select()
